$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column S: 2022 data, added to the right of the existing 2021 (R) column ---

# S4 header (2022) - clone formatting from R4 (same header row style)
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S4").Value = 2022

# S5 (bold row, font 7 family) - build the new number-format style explicitly
$cell = $ws.Range("S5")
$cell.NumberFormat = "0.0"
$cell.Font.Name = "Times New Roman"
$cell.Font.Size = 9
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4152
$cell.VerticalAlignment = -4108
$cell.Value = 4.9000000000000004

# S6 (regular row, font 11 family) - new number-format style explicitly
$cell = $ws.Range("S6")
$cell.NumberFormat = "0.0"
$cell.Font.Name = "Times New Roman"
$cell.Font.Size = 9
$cell.Font.Bold = $false
$cell.HorizontalAlignment = -4152
$cell.VerticalAlignment = -4108
$cell.Value = 3.4

# S7:S13 share S6's exact new style - clone it via copy/paste-special so no
# extra intermediate styles get minted, then fill in each value.
$ws.Range("S6").Copy() | Out-Null
$ws.Range("S7:S13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("S7").Value = 3.5
$ws.Range("S8").Value = 13.1
$ws.Range("S9").Value = 8.1
$ws.Range("S10").Value = 2.5
$ws.Range("S11").Value = 2.6
$ws.Range("S12").Value = 10.8
$ws.Range("S13").Value = 2.1

# S14 (bottom, thick-bottom-border row) - clone formatting from R14
$ws.Range("R14").Copy() | Out-Null
$ws.Range("S14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S14").Value = 1.1000000000000001

# Move the active selection the way the author's session ended up (T4 instead of T9)
$ws.Range("T4").Select() | Out-Null
